$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -0
$ws.Range("B2").Value = -0.08002405806976164
$ws.Range("D2").Value = 0.230202729308465
$ws.Range("E2").Value = 0.005108318808481852
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = -0
$ws.Range("J2").Value = -0
$ws.Range("K2").Value = -0.03860066713380397
$ws.Range("L2").Value = -0
$ws.Range("M2").Value = 0.2075828185783251
$ws.Range("N2").Value = 0.01204929019375068
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = -0.09806669489094141
$ws.Range("V2").Value = 0.01229636434604111
$ws.Range("W2").Value = -0.004991894583601985
$ws.Range("Z2").Value = -0
$ws.Range("AA2").Value = -0
$ws.Range("AB2").Value = 0
$ws.Range("AC2").Value = 0.03388509010844856
$ws.Range("AD2").Value = 0
$ws.Range("AE2").Value = 0.006056585975831554
$ws.Range("AF2").Value = -0.004450682315163321
$ws.Range("AG2").Value = -0
$ws.Range("AI2").Value = -0
$ws.Range("AJ2").Value = 0
$ws.Range("AK2").Value = 0
$ws.Range("AL2").Value = -0.03033852235949213
$ws.Range("AM2").Value = 0
$ws.Range("AN2").Value = 0.03456439921250479
$ws.Range("AO2").Value = 0.06390532446227024
$ws.Range("AQ2").Value = 0
$ws.Range("AR2").Value = -0
$ws.Range("AS2").Value = -0
$ws.Range("AT2").Value = 0
$ws.Range("AU2").Value = -0.1606877030427201
$ws.Range("AW2").Value = 0.07871500957670378
$ws.Range("AX2").Value = 0.005945696286255524
$ws.Range("AY2").Value = -0
$ws.Range("BB2").Value = -0
$ws.Range("BC2").Value = -0
$ws.Range("BD2").Value = -0.01757312437285266
$ws.Range("BF2").Value = 0.07897551879926892
$ws.Range("BG2").Value = 0.02712815448214431
$ws.Range("BJ2").Value = -0
$ws.Range("BL2").Value = 0
$ws.Range("BM2").Value = 0.03040979083744841
$ws.Range("BO2").Value = -0.04594451529409452
$ws.Range("BP2").Value = -0.08000625774407112
$ws.Range("BR2").Value = -0
$ws.Range("BU2").Value = 0
$ws.Range("BV2").Value = -0.06713017421494699
$ws.Range("BW2").Value = 0
$ws.Range("BX2").Value = 0.02182646807514023
$ws.Range("BY2").Value = -0.01884615699057904
$ws.Range("BZ2").Value = -0
$ws.Range("CB2").Value = 0
$ws.Range("CD2").Value = -0
$ws.Range("CE2").Value = 0.03283043049272813
$ws.Range("CF2").Value = 0
$ws.Range("CG2").Value = -0.03616967539422405
$ws.Range("CH2").Value = 0.01616340536772806
$ws.Range("CJ2").Value = -0
$ws.Range("CL2").Value = 0
$ws.Range("CM2").Value = -0
$ws.Range("CN2").Value = -0.007863134234805585
$ws.Range("CO2").Value = -0
$ws.Range("CP2").Value = 0.02819151524651156
$ws.Range("CQ2").Value = 0.02816355631835672
$ws.Range("CT2").Value = 0
$ws.Range("CU2").Value = -0
$ws.Range("CV2").Value = -0
$ws.Range("CW2").Value = 0.04648971561772923
$ws.Range("CX2").Value = -0
$ws.Range("CY2").Value = -0.04027009503825581
$ws.Range("CZ2").Value = 0.01093820471416725
$ws.Range("DE2").Value = -0
$ws.Range("DF2").Value = 0.03395264743488322
$ws.Range("DH2").Value = 0.01468105289247713
$ws.Range("DI2").Value = 0.03297448683683691
$ws.Range("DJ2").Value = 0
$ws.Range("DK2").Value = -0
$ws.Range("DL2").Value = -0
$ws.Range("DN2").Value = 0
$ws.Range("DO2").Value = -0.01918322500336656
$ws.Range("DQ2").Value = 0.04349730583284731
$ws.Range("DR2").Value = -0.005621113725438998
$ws.Range("DS2").Value = -0
$ws.Range("DT2").Value = 0
$ws.Range("DU2").Value = -0
$ws.Range("DV2").Value = -0
$ws.Range("DW2").Value = 0
$ws.Range("DX2").Value = -0.06400284214897392
$ws.Range("DY2").Value = -0
$ws.Range("DZ2").Value = -0.0109657149773436
$ws.Range("EA2").Value = -0.02088623205835872
$ws.Range("EB2").Value = 0
$ws.Range("ED2").Value = 0
$ws.Range("EF2").Value = -0
$ws.Range("EG2").Value = 0.03911788968558339
$ws.Range("EH2").Value = 0
$ws.Range("EI2").Value = 0.08937657361646349
$ws.Range("EJ2").Value = -0.03570638022466691
$ws.Range("EO2").Value = 0
$ws.Range("EP2").Value = 0.04835681117410522
$ws.Range("EQ2").Value = 0
$ws.Range("ER2").Value = -0.0470410529826082
$ws.Range("ES2").Value = 0.01999236228472831
$ws.Range("ET2").Value = 0
$ws.Range("EV2").Value = 0
$ws.Range("EX2").Value = 0
$ws.Range("EY2").Value = 0.04098487980350489
$ws.Range("EZ2").Value = 0
$ws.Range("FA2").Value = -0.03239473819547633
$ws.Range("FB2").Value = 0.01641101982812671
$ws.Range("FD2").Value = -0
$ws.Range("FG2").Value = -0
$ws.Range("FH2").Value = -0.01242405373375696
$ws.Range("FJ2").Value = -0.00543148769200911
$ws.Range("FK2").Value = 0.01687241999401279
$ws.Range("FL2").Value = -0
$ws.Range("FM2").Value = 0
$ws.Range("FP2").Value = -0
$ws.Range("FQ2").Value = -0.01297647463482977
$ws.Range("FR2").Value = -0
$ws.Range("FS2").Value = 0.006666584027607884
$ws.Range("FT2").Value = -0.008699428982887088
$ws.Range("FV2").Value = -0
$ws.Range("FW2").Value = -0
$ws.Range("FY2").Value = 0
$ws.Range("FZ2").Value = -0.03139313018916433
$ws.Range("GB2").Value = 0.01992368484076844
$ws.Range("GD2").Value = 0
$ws.Range("GE2").Value = -0
